$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2-12: I = 1, J = copy of H column value for that row
for ($r = 2; $r -le 12; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
